$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C ("날짜"), shifting the old C ("구분") and D ("메뉴") right by one.
$ws.Columns.Item(3).Insert()

# Insert a new column at F ("칼로리"), after the (now) E ("메뉴") column.
$ws.Columns.Item(6).Insert()

# Header row
$ws.Range("C1").Value = "날짜"
$ws.Range("F1").Value = "칼로리"

# New row 6 index + style to match the other index cells (A2:A5)
$ws.Range("A6").Value = 4
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

# Row 2
$ws.Range("B2").Value = "휘봉고등학교"
$ws.Range("D2").Value = "중식"
$ws.Range("E2").Value = "현미밥, 청양콩나물국, 청포묵무침, 제육볶음, 배추김치, 귤"
$ws.Range("F2").Value = "356.7 Kcal"

# Row 3
$ws.Range("B3").Value = "휘경여자고등학교"
$ws.Range("D3").Value = "중식"
$ws.Range("E3").Value = "옥수수밥, 쇠고기샤브샤브국, 진미채고추장볶음, 치즈불닭, 굴림만두, 깍두기, 자몽데르뜨"
$ws.Range("F3").Value = "801.3 Kcal"

# Row 4
$ws.Range("B4").Value = "휘경여자고등학교"
$ws.Range("D4").Value = "중식"
$ws.Range("E4").Value = "흑미밥, 사골조랭이고기만두국, 도토리묵야채무침, 스팸감자구이, 어향동태강정, 배추김치"
$ws.Range("F4").Value = "935.2 Kcal"

# Row 5
$ws.Range("B5").Value = "휘경공업고등학교"
$ws.Range("D5").Value = "중식"
$ws.Range("E5").Value = "칼슘강화강낭콩밥, 맑은콩나물국, 시금치나물무침, 닭볶음탕, 어묵피망볶음, 깍두기"
$ws.Range("F5").Value = "1022.3 Kcal"

# New row 6
$ws.Range("B6").Value = "휘경공업고등학교"
$ws.Range("D6").Value = "중식"
$ws.Range("E6").Value = "칼슘강화현미밥, 아욱국, 돈등뼈김치찜, 갈릭난*커리소스, 깍두기, 리코타치즈샐러드"
$ws.Range("F6").Value = "1099.0 Kcal"

# Dates (column C): these look like plain numbers, so force Text format
# first - otherwise Excel auto-converts the digit string to a number -
# then strip the format override again so the cells keep the plain
# (unstyled) look of their neighbours, just stored as text.
$ws.Range("C2:C6").NumberFormat = "@"
$ws.Range("C2").Value = "20210104"
$ws.Range("C3").Value = "20210104"
$ws.Range("C4").Value = "20210105"
$ws.Range("C5").Value = "20210104"
$ws.Range("C6").Value = "20210105"
$ws.Range("C2:C6").ClearFormats()
